# correction on 3 & 4 stages
# G column (count 3) formulas keep "Count_Number":"3" but the "id" offset
# moves from ROW()+6 to ROW()+3.
# K column formulas change "Count_Number" from "5" to "4" and the "id"
# offset moves from ROW()+12/ROW()+6 to ROW()+9.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9 - standalone (non-shared) formulas
$ws.Range("G9").Formula = '="{""Candidate_First_Pref_Votes"":"""&$D1&""",""Status"":"""",""Occurred_On_Count"":"""",""Surname"":"""&$B1&""",""Firstname"":"""&$A1&""",""Constituency_Number"":""2"",""Party_Name"":"""&$C1&""",""Candidate_Id"":"""&ROW()&""",""Count_Number"":""3"",""Transfers"":"""&G1&""",""id"":"&ROW()+3&",""Total_Votes"":"""&$H1&"""},"'
$ws.Range("K9").Formula = '="{""Candidate_First_Pref_Votes"":"""&$D1&""",""Status"":"""",""Occurred_On_Count"":"""",""Surname"":"""&$B1&""",""Firstname"":"""&$A1&""",""Constituency_Number"":""2"",""Party_Name"":"""&$C1&""",""Candidate_Id"":"""&ROW()&""",""Count_Number"":""4"",""Transfers"":"""&K1&""",""id"":"&ROW()+9&",""Total_Votes"":"""&L1&"""},"'

# Rows 10-14 - shared formula groups (fill the whole range so the shared
# formula block si=2 / si=3 is rewritten consistently across all rows)
$ws.Range("G10:G14").Formula = '="{""Candidate_First_Pref_Votes"":"""&$D2&""",""Status"":"""",""Occurred_On_Count"":"""",""Surname"":"""&$B2&""",""Firstname"":"""&$A2&""",""Constituency_Number"":""2"",""Party_Name"":"""&$C2&""",""Candidate_Id"":"""&ROW()&""",""Count_Number"":""3"",""Transfers"":"""&G2&""",""id"":"&ROW()+3&",""Total_Votes"":"""&$H2&"""},"'
$ws.Range("K10:K14").Formula = '="{""Candidate_First_Pref_Votes"":"""&$D2&""",""Status"":"""",""Occurred_On_Count"":"""",""Surname"":"""&$B2&""",""Firstname"":"""&$A2&""",""Constituency_Number"":""2"",""Party_Name"":"""&$C2&""",""Candidate_Id"":"""&ROW()&""",""Count_Number"":""4"",""Transfers"":"""&K2&""",""id"":"&ROW()+9&",""Total_Votes"":"""&L2&"""},"'
